# Coachingslijst.xlsx update: a coaching record (P-nr 33085, Van Hoe Ivo)
# that was completed got moved from the active "Coaching" list to the
# "Voltooide coachings" (completed) list.

$wb = $excel.ActiveWorkbook

$wsCoaching = $wb.Worksheets.Item("Coaching")
$wsVoltooid = $wb.Worksheets.Item("Voltooide coachings")

# 1) Remove the now-completed entry (row 74, P-nr 33085 "Van Hoe Ivo") from
#    the "Coaching" sheet. This shifts every row below it up by one and
#    Excel automatically shrinks the Tabel2 structured range / autofilter /
#    sheet dimension to match.
$wsCoaching.Rows("74:74").Delete()

# 2) Record that same coaching as completed on the "Voltooide coachings"
#    sheet (row 197 there was a blank template row waiting for the next
#    completed entry - P-nr drives the VLOOKUP formulas already in place).
$wsVoltooid.Range("B197").Value = 33085
$wsVoltooid.Range("I197").Value = "schadegevallen"
$wsVoltooid.Range("J197").Value = "bus"
$wsVoltooid.Range("K197").Value = "8538 Philippe Vandewalle"
$wsVoltooid.Range("L197").Value = 45950
$wsVoltooid.Range("M197").Value = "wegcode"
$wsVoltooid.Range("O197").Value = "goed"

# 3) Update on-screen navigation state: the user ended up scrolled further
#    down the "Coaching" sheet, and "Voltooide coachings" became the active
#    (visible) tab with its own scrolled selection.
$wsCoaching.Range("A74:XFD74").Select() | Out-Null

$wsVoltooid.Activate() | Out-Null
$wsVoltooid.Range("A204:XFD204").Select() | Out-Null
